# Generate Report for Handoff
# Updates the handoff report: Priority flips from "low" to "ht" for the
# rows that have now been handed off (0d44ddd4 / 42069519 / 4ad57348 /
# b4a4978e), and the handoff/generation timestamps advance accordingly
# on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) - shared
# across rows 4-7, advance the generation timestamp.
$overview.Range("G4").Value = "2016-08-15 18:28:37"
$overview.Range("G5").Value = "2016-08-15 18:28:37"
$overview.Range("G6").Value = "2016-08-15 18:28:37"
$overview.Range("G7").Value = "2016-08-15 18:28:37"

# zh-cn / de-de sheets: Priority column (E) for rows 4-7 flips from
# "low" to "ht" now that the files are ready for handoff.
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("E5").Value = "ht"
$zhcn.Range("E6").Value = "ht"
$zhcn.Range("E7").Value = "ht"

$dede.Range("E4").Value = "ht"
$dede.Range("E5").Value = "ht"
$dede.Range("E6").Value = "ht"
$dede.Range("E7").Value = "ht"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for rows 4-7 -
# shares the same underlying handoff timestamp, advances together.
$zhcn.Range("H4").Value = "2016-08-15 18:28:32"
$zhcn.Range("H5").Value = "2016-08-15 18:28:32"
$zhcn.Range("H6").Value = "2016-08-15 18:28:32"
$zhcn.Range("H7").Value = "2016-08-15 18:28:32"

# de-de sheet: "Latest Handoff Datetime" column (H) for rows 4-7 -
# shares the same underlying generation timestamp as the Overview sheet.
$dede.Range("H4").Value = "2016-08-15 18:28:37"
$dede.Range("H5").Value = "2016-08-15 18:28:37"
$dede.Range("H6").Value = "2016-08-15 18:28:37"
$dede.Range("H7").Value = "2016-08-15 18:28:37"
